$d = $word.ActiveDocument

# The first two paragraphs of the document are:
#   1) Heading1 "Day After Day - May 1941" (wrapped in a bookmark)
#   2) Bold "By Dorothy Day"
# They need to become a pandoc-style title block:
#   1) Title-styled paragraph "Day After Day - May 1941" split word-by-word
#      into separate runs, bookmark removed
#   2) Authors-styled paragraph "Dorothy Day" split word-by-word into
#      separate runs, bold removed, "By " prefix removed

$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)

$titleWords = @("Day", " ", "After", " ", "Day", " ", "-", " ", "May", " ", "1941")
$authorWords = @("Dorothy", " ", "Day")

$titleRuns = ""
foreach ($w in $titleWords) {
    $titleRuns += '<w:r><w:t xml:space="preserve">' + $w + '</w:t></w:r>'
}

$authorRuns = ""
foreach ($w in $authorWords) {
    $authorRuns += '<w:r><w:t xml:space="preserve">' + $w + '</w:t></w:r>'
}

$bodyFragment = '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $titleRuns + '</w:p>' + `
                '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + $authorRuns + '</w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $bodyFragment + '</w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($packageXml) | Out-Null
